$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.529.15'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '2.486.04'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '315.13'
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = '93.80'
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("D7").Value = '0.543'
$ws.Range("E7").Value = '  -0.92%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = '0.507'
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").Value = '32.91'
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").Value = '0.0788'
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("D12").Value = '0.111'
$ws.Range("E12").Value = '  +2.94%  '
$ws.Range("D13").Value = '2.868.57'
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '6.86'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").Value = '16.15'
$ws.Range("E15").Value = '  +11.11%  '
$ws.Range("D16").Value = '2.457.44'
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").Value = '0.762'
$ws.Range("E17").Value = '  -1.78%  '
$ws.Range("D18").Value = '41.555.36'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0942'
$ws.Range("E19").Value = '  +2.94%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '6.38'
$ws.Range("E20").Value = '  +1.46%  '
$ws.Range("D21").Value = '71.28'
$ws.Range("E21").Value = '  +4.86%  '
$ws.Range("D22").Value = '11.45'
$ws.Range("E22").Value = '  +2.16%  '
$ws.Range("D23").Value = '237.83'
$ws.Range("E23").Value = '  +1.08%  '
$ws.Range("D24").Value = '2.73'
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = '1.91'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '25.04'
$ws.Range("E27").Value = '  +5.21%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").Value = '9.72'
$ws.Range("E29").Value = '  +1.62%  '
$ws.Range("D30").Value = '36.30'
$ws.Range("E30").Value = '  +2.29%  '
$ws.Range("D31").Value = '157.23'
$ws.Range("E31").Value = '  +3.39%  '
$ws.Range("D32").Value = '5.49'
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").Value = '2.57'
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("D34").Value = '0.0756'
$ws.Range("E34").Value = '  +3.00%  '
$ws.Range("D35").Value = '17.80'
$ws.Range("E35").Value = '  +5.05%  '
$ws.Range("E36").Value = '  -6.25%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '2.94'
$ws.Range("E37").Value = '  -1.48%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.106'
$ws.Range("E38").Value = '  +3.57%  '
$ws.Range("D39").Value = '1.85'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").Value = '4.12'
$ws.Range("E41").Value = '  -2.22%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").Value = '19.99'
$ws.Range("E43").Value = '  -1.05%  '
$ws.Range("D44").Value = '1.973.19'
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").Value = '0.0286'
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("D46").Value = '2.98'
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("D47").Value = '8.99'
$ws.Range("E47").Value = '  +4.44%  '
$ws.Range("D48").Value = '2.723.53'
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("D49").Value = '97.63'
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").Value = '68.85'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("D51").Value = '72.96'
$ws.Range("E51").Value = '  -1.40%  '
